$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the D:E data range so price/volume strings that look
# like numbers ("21.20", "0.820", "7.30", ...) are not silently reinterpreted
# by Excel as numeric values (which would drop trailing zeros / precision).
$rngD = $ws.Range("D2:D51")
$rngE = $ws.Range("E2:E51")
$rngD.NumberFormat = "@"
$rngE.NumberFormat = "@"

$ws.Range("D2").Value = "59.771.33"
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").Value = "2.647.30"
$ws.Range("E3").Value = "  +2.06%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "536.34"
$ws.Range("E5").Value = "  +1.06%  "
$ws.Range("D6").Value = "145.05"
$ws.Range("E6").Value = "  +3.40%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +0.99%  "
$ws.Range("D9").Value = "2.663.59"
$ws.Range("E9").Value = "  +2.18%  "
$ws.Range("D10").Value = "6.71"
$ws.Range("E10").Value = "  +4.09%  "
$ws.Range("E11").Value = "  +1.53%  "
$ws.Range("E12").Value = "  +1.32%  "
$ws.Range("E13").Value = "  -1.24%  "
$ws.Range("D14").Value = "3.119.77"
$ws.Range("E14").Value = "  +2.19%  "
$ws.Range("D15").Value = "59.688.82"
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("D16").Value = "21.20"
$ws.Range("E16").Value = "  +3.52%  "
$ws.Range("D17").Value = "2.636.19"
$ws.Range("E17").Value = "  +2.03%  "
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("D19").Value = "344.19"
$ws.Range("E19").Value = "  -0.81%  "
$ws.Range("E20").Value = "  +1.88%  "
$ws.Range("D21").Value = "10.25"
$ws.Range("E21").Value = "  +1.34%  "
$ws.Range("D22").Value = "6.35"
$ws.Range("E22").Value = "  -0.85%  "
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "66.86"
$ws.Range("E24").Value = "  -0.96%  "
$ws.Range("E25").Value = "  +2.27%  "
$ws.Range("E26").Value = "  -0.79%  "
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("D28").Value = "7.30"
$ws.Range("E28").Value = "  +2.00%  "
$ws.Range("D29").Value = "0.0₃0748"
$ws.Range("E29").Value = "  +1.42%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  +2.66%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "19.07"
$ws.Range("E32").Value = "  +1.50%  "
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").Value = "5.84"
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").Value = "150.21"
$ws.Range("E34").Value = "  +0.96%  "
$ws.Range("D35").Value = "4.04"
$ws.Range("E35").Value = "  +1.38%  "
$ws.Range("D36").Value = "1.15"
$ws.Range("E36").Value = "  +3.11%  "
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("D38").Value = "0.838"
$ws.Range("E38").Value = "  +0.66%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "294.62"
$ws.Range("E39").Value = "  +8.71%  "
$ws.Range("B40").Value = "SuiNetwork"
$ws.Range("C40").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D40").Value = "0.820"
$ws.Range("E40").Value = "  -0.71%  "
$ws.Range("D41").Value = "3.59"
$ws.Range("E41").Value = "  +1.81%  "
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").Value = "0.604"
$ws.Range("E43").Value = "  +1.42%  "
$ws.Range("D44").Value = "0.0545"
$ws.Range("E44").Value = "  +5.02%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "19.38"
$ws.Range("E45").Value = "  +5.04%  "
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").Value = "10.73"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("E47").Value = "  -0.44%  "
$ws.Range("E48").Value = "  +2.28%  "
$ws.Range("D49").Value = "1.970.72"
$ws.Range("E49").Value = "  +0.92%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "4.56"
$ws.Range("E50").Value = "  -1.83%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "18.37"
$ws.Range("E51").Value = "  +0.64%  "

# Restore default (Normal) style on the range now that values are written,
# so cells do not keep an explicit style index pointing at the Text format.
$rngD.Style = "Normal"
$rngE.Style = "Normal"
